$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.9
$ws.Range("I2").Value = 2.45
$ws.Range("J2").Value = 3.5
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = 1.22
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 1.91
$ws.Range("Y2").Value = 8
$ws.Range("AC2").Value = 26
$ws.Range("AD2").Value = 41
$ws.Range("AG2").Value = 17
$ws.Range("AI2").Value = 700
$ws.Range("AK2").Value = 11

# Row 3 updates
$ws.Range("G3").Value = 2.63
$ws.Range("I3").Value = 2.88
$ws.Range("AP3").Value = 1.9
$ws.Range("AQ3").Value = 1.95
